# Remove obsolete finite-difference translation rows:
#   strRadBackwardOne, strRadCentralFive, strRadCentralThree, strRadForwardOne
# These currently sit at worksheet rows 127-129 and 131 (row 130,
# strRadCurrentCulture, is kept). Delete bottom-up so earlier row numbers
# stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Delete()
$ws.Rows.Item(129).Delete()
$ws.Rows.Item(128).Delete()
$ws.Rows.Item(127).Delete()
